$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3407069999999999
$ws.Range("N2").Value = 1.022121
$ws.Range("O2").Value = 0.1055965976712818
$ws.Range("P2").Value = 0.1055965976712818
$ws.Range("Q2").Value = 14.566181182394
$ws.Range("R2").Value = 131.095630641546
$ws.Range("S2").Value = 0.09006635879467877
$ws.Range("T2").Value = 0.09006635879467877
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.3782153560188308
$ws.Range("P3").Value = 0.3782153560188308
$ws.Range("Q3").Value = 52.17169419495621
$ws.Range("R3").Value = 469.5452477546059
$ws.Range("S3").Value = 0.3225906961783997
$ws.Range("T3").Value = 0.3225906961783998
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5161880463098875
$ws.Range("P4").Value = 0.5161880463098875
$ws.Range("Q4").Value = 71.20389077441509
$ws.Range("R4").Value = 640.8350169697359
$ws.Range("S4").Value = 0.440271550501995
$ws.Range("T4").Value = 0.440271550501995
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3407069999999999
$ws.Range("N5").Value = 1.022121
$ws.Range("O5").Value = 0.1055965976712818
$ws.Range("P5").Value = 0.1055965976712818
$ws.Range("Q5").Value = 0.7928574425959999
$ws.Range("R5").Value = 7.135716983364
$ws.Range("S5").Value = 0.0049024368160541
$ws.Range("T5").Value = 0.004902436816054099
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.3782153560188308
$ws.Range("P6").Value = 0.3782153560188308
$ws.Range("Q6").Value = 2.839777668378222
$ws.Range("S6").Value = 0.017559058971916
$ws.Range("T6").Value = 0.017559058971916
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5161880463098875
$ws.Range("P7").Value = 0.5161880463098875
$ws.Range("Q7").Value = 3.875726522647111
$ws.Range("R7").Value = 34.881538703824
$ws.Range("S7").Value = 0.02396459107626013
$ws.Range("T7").Value = 0.02396459107626013
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3407069999999999
$ws.Range("N8").Value = 1.022121
$ws.Range("O8").Value = 0.1055965976712818
$ws.Range("P8").Value = 0.1055965976712818
$ws.Range("Q8").Value = 1.718804806325999
$ws.Range("R8").Value = 15.469243256934
$ws.Range("S8").Value = 0.01062780206054892
$ws.Range("T8").Value = 0.01062780206054892
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.3782153560188308
$ws.Range("P9").Value = 0.3782153560188308
$ws.Range("Q9").Value = 6.156243535185999
$ws.Range("R9").Value = 55.40619181667399
$ws.Range("S9").Value = 0.038065600868515
$ws.Range("T9").Value = 0.038065600868515
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5161880463098875
$ws.Range("P10").Value = 0.5161880463098875
$ws.Range("Q10").Value = 8.402036756215997
$ws.Range("R10").Value = 75.61833080594398
$ws.Range("S10").Value = 0.05195190473163237
$ws.Range("T10").Value = 0.05195190473163237
